# Applies the "Automatic update of files." edit:
#   1. Updates the "Förändrad" (changed) date in column C for every data
#      row (rows 2-66) from 45184 to 45186.
#   2. Adds a friendly display-text second argument to the HYPERLINK()
#      formulas on row 2 (columns S, T, V, W, X, Y) using the
#      designation in column A of that row ("A 30186-2021").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column C: "Förändrad" date, rows 2-66 -----------------------------
for ($row = 2; $row -le 66; $row++) {
    $ws.Cells.Item($row, 3).Value = 45186
}

# --- 2. Row 2 hyperlink formulas: add display text ------------------------
$designation = $ws.Cells.Item(2, 1).Value2

$hyperlinkCols = 19, 20, 22, 23, 24, 25   # S, T, V, W, X, Y

foreach ($col in $hyperlinkCols) {
    $cell = $ws.Cells.Item(2, $col)
    $formula = $cell.Formula
    if ($formula -match '^=HYPERLINK\("([^"]*)"\)$') {
        $url = $Matches[1]
        $cell.Formula = '=HYPERLINK("' + $url + '", "' + $designation + '")'
    }
}
